$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "(August 2009- May 2013)" -> "(August 2009 - May 2013)"
#    (a missing space before the en-dash is added). The canonical edit splits
#    the surrounding run into three runs; we replicate that by inserting the
#    space and then nudging a (no-op) character-formatting property on the
#    newly created pieces so Word materialises them as separate <w:r> runs
#    instead of silently re-merging them into the original run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("(August 2009- May 2013)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $rng.Start
$end = $rng.End

$insertPoint = $d.Range($start + 12, $start + 12)
$insertPoint.InsertAfter(" ")

$spaceRun = $d.Range($start + 12, $start + 13)
$spaceRun.Font.Bold = 1
$spaceRun.Font.Bold = 0

$tailRun = $d.Range($start + 13, $end + 1)
$tailRun.Font.Bold = 1
$tailRun.Font.Bold = 0

# ---------------------------------------------------------------------------
# 2) New bullet "Full-Stack Web Developer" added under the Qualifications
#    list, right after "Can type 70 wpm".
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Can type 70 wpm", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $rng.Paragraphs(1)
$anchorPara.Range.InsertParagraphAfter()
$newPara = $anchorPara.Next()
$newPara.Range.Text = "Full-Stack Web Developer"

# ---------------------------------------------------------------------------
# 3) "Technical Skills:" paragraph gains a trailing space run (bold), merging
#    what used to be a separate empty bold paragraph into this one.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Technical Skills:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$end = $rng.End
$insertPoint = $d.Range($end, $end)
$insertPoint.InsertAfter(" ")

$spaceRun = $d.Range($end, $end + 1)
$spaceRun.Font.Bold = 0
$spaceRun.Font.Bold = 1

Write-Output "done"
